$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.144.98'
$ws.Range("D3").Value = '1.670.11'
$ws.Range("E3").Value = '  -1.50%  '
$ws.Range("E4").Value = '  -0.59%  '
$ws.Range("E6").Value = '  -2.89%  '
$ws.Range("E8").Value = '  -3.63%  '
$ws.Range("E9").Value = '  -2.25%  '
$ws.Range("E10").Value = '  -2.02%  '
$ws.Range("D11").Value = '''0.07569'
$ws.Range("E11").Value = '  -1.29%  '
$ws.Range("D12").Value = '1.682.05'
$ws.Range("E12").Value = '  -1.39%  '
$ws.Range("E13").Value = '  -2.04%  '
$ws.Range("E14").Value = '  -4.03%  '
$ws.Range("E15").Value = '  -0.17%  '
$ws.Range("E16").Value = '  -5.35%  '
$ws.Range("D17").Value = '26.184.79'
$ws.Range("E17").Value = '  -0.89%  '
$ws.Range("E19").Value = '  -3.52%  '
$ws.Range("D20").Value = '''187.17'
$ws.Range("E20").Value = '  -2.03%  '
$ws.Range("E21").Value = '  -4.76%  '
$ws.Range("E22").Value = '  -1.55%  '
$ws.Range("E23").Value = '  -0.72%  '
$ws.Range("E24").Value = '  +0.24%  '
$ws.Range("E25").Value = '  -2.81%  '
$ws.Range("E26").Value = '  -4.05%  '
$ws.Range("D27").Value = '''16.00'
$ws.Range("E27").Value = '  +0.58%  '
$ws.Range("D28").Value = '''0.06282'
$ws.Range("E28").Value = '  -1.12%  '
$ws.Range("E29").Value = '  -2.30%  '
$ws.Range("E30").Value = '  -3.44%  '
$ws.Range("D31").Value = '''3.511'
$ws.Range("E31").Value = '  -2.89%  '
$ws.Range("E32").Value = '  -5.00%  '
$ws.Range("E33").Value = '  -3.66%  '
$ws.Range("E35").Value = '  -2.28%  '
$ws.Range("E36").Value = '  -0.41%  '
$ws.Range("E37").Value = '  -0.58%  '
$ws.Range("E38").Value = '  +0.25%  '
$ws.Range("D39").Value = '1.104.04'
$ws.Range("E39").Value = '  -1.09%  '
$ws.Range("D40").Value = '''0.01612'
$ws.Range("E40").Value = '  -2.53%  '
$ws.Range("D41").Value = '''0.8767'
$ws.Range("E41").Value = '  -1.14%  '
$ws.Range("E42").Value = '  -1.04%  '
$ws.Range("E43").Value = '  -0.95%  '
$ws.Range("D44").Value = '1.822.65'
$ws.Range("E44").Value = '  -1.31%  '
$ws.Range("E45").Value = '  +0.24%  '
$ws.Range("E47").Value = '  -0.36%  '
$ws.Range("D48").Value = '''8.034'
$ws.Range("E48").Value = '  -1.86%  '
$ws.Range("E49").Value = '  -0.99%  '
$ws.Range("D50").Value = '''0.4248'
$ws.Range("E50").Value = '  -1.25%  '
$ws.Range("E51").Value = '  -1.75%  '
